$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.688.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.611.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "631.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.610.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.58%  "
$ws.Range("E13").Value = "  +6.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.223.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.614.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.542.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +6.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.22%  "
$ws.Range("E23").Value = "  +5.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("E25").Value = "  +11.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.757.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.96%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.16%  "
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.173"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.15%  "
$ws.Range("E33").Value = "  +9.24%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +7.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.609.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.13%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "179.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +6.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +24.32%  "
$ws.Range("E46").Value = "  +5.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.03%  "
$ws.Range("E48").Value = "  +13.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("E51").Value = "  +11.98%  "
